$d = $word.ActiveDocument

$d.Content.Find.Execute("A fórmula para calulcar o MAE é:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A fórmula para calcular o MAE é:", 2)
